# Generate Report for Handoff
# Adds two new source files (285cd6e1-... and 7dc73b0e-...) to the
# localization-status workbook: one new row per file on the "Overview"
# sheet, and one new row per file on each of the language sheets
# ("zh-cn", "de-de").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Shared data for the two newly handed-off files
# ---------------------------------------------------------------------
$file1Md      = "285cd6e1-5894-4e53-aa44-4e460617a46b.md"
$file1ZhXlf   = "285cd6e1-5894-4e53-aa44-4e460617a46b.d4f36eb6fb25f6a5020d6d9f34756ca09dc953c5.zh-cn.xlf"
$file1DeXlf   = "285cd6e1-5894-4e53-aa44-4e460617a46b.d4f36eb6fb25f6a5020d6d9f34756ca09dc953c5.de-de.xlf"

$file2Md      = "7dc73b0e-ad13-4087-b4f7-51b36b638f20.md"
$file2ZhXlf   = "7dc73b0e-ad13-4087-b4f7-51b36b638f20.119f32a2eadec3626a1ddae8ad7561e4ff39bef0.zh-cn.xlf"
$file2DeXlf   = "7dc73b0e-ad13-4087-b4f7-51b36b638f20.119f32a2eadec3626a1ddae8ad7561e4ff39bef0.de-de.xlf"

$handoffDateTime   = "2016-03-25 07:22:34"
$handbackDateTimeZh = "2016-03-25 07:22:30"
$handbackDateTimeDe = "2016-03-25 07:22:34"
$noHandbackYet      = "0001-01-01 00:00:00"

$status   = "Ready for handoff"
$ext      = ".md"
$reason   = "Include"
$dateFmt  = "yyyy-mm-dd HH:mm:ss"

$ghBase1 = "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e"
$ghZhHandoffBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$ghDeHandoffBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

# ---------------------------------------------------------------------
# Sheet "Overview" -- one row per file, columns A:D
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status
$wsOverview.Range("D4").Value = $handoffDateTime
$wsOverview.Range("D4").NumberFormat = $dateFmt
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "$ghBase1/$file1Md", [Type]::Missing, [Type]::Missing, $file1Md) | Out-Null

$wsOverview.Range("B5").Value = $status
$wsOverview.Range("C5").Value = $status
$wsOverview.Range("D5").Value = $handoffDateTime
$wsOverview.Range("D5").NumberFormat = $dateFmt
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "$ghBase1/$file2Md", [Type]::Missing, [Type]::Missing, $file2Md) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" -- one row per file, columns A:L
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B4").Value = $ext
$wsZh.Range("C4").Value = $status
$wsZh.Range("E4").Value = $handbackDateTimeZh
$wsZh.Range("E4").NumberFormat = $dateFmt
$wsZh.Range("H4").Value = $noHandbackYet
$wsZh.Range("H4").NumberFormat = $dateFmt
$wsZh.Range("J4").Value = $reason
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "$ghBase1/$file1Md", [Type]::Missing, [Type]::Missing, $file1Md) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), "$ghZhHandoffBase/$file1ZhXlf", [Type]::Missing, [Type]::Missing, $file1ZhXlf) | Out-Null

$wsZh.Range("B5").Value = $ext
$wsZh.Range("C5").Value = $status
$wsZh.Range("E5").Value = $handbackDateTimeZh
$wsZh.Range("E5").NumberFormat = $dateFmt
$wsZh.Range("H5").Value = $noHandbackYet
$wsZh.Range("H5").NumberFormat = $dateFmt
$wsZh.Range("J5").Value = $reason
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "$ghBase1/$file2Md", [Type]::Missing, [Type]::Missing, $file2Md) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), "$ghZhHandoffBase/$file2ZhXlf", [Type]::Missing, [Type]::Missing, $file2ZhXlf) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" -- one row per file, columns A:L
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B4").Value = $ext
$wsDe.Range("C4").Value = $status
$wsDe.Range("E4").Value = $handbackDateTimeDe
$wsDe.Range("E4").NumberFormat = $dateFmt
$wsDe.Range("H4").Value = $noHandbackYet
$wsDe.Range("H4").NumberFormat = $dateFmt
$wsDe.Range("J4").Value = $reason
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "$ghBase1/$file1Md", [Type]::Missing, [Type]::Missing, $file1Md) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), "$ghDeHandoffBase/$file1DeXlf", [Type]::Missing, [Type]::Missing, $file1DeXlf) | Out-Null

$wsDe.Range("B5").Value = $ext
$wsDe.Range("C5").Value = $status
$wsDe.Range("E5").Value = $handbackDateTimeDe
$wsDe.Range("E5").NumberFormat = $dateFmt
$wsDe.Range("H5").Value = $noHandbackYet
$wsDe.Range("H5").NumberFormat = $dateFmt
$wsDe.Range("J5").Value = $reason
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "$ghBase1/$file2Md", [Type]::Missing, [Type]::Missing, $file2Md) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), "$ghDeHandoffBase/$file2DeXlf", [Type]::Missing, [Type]::Missing, $file2DeXlf) | Out-Null

Write-Host "Generated handoff report rows for $file1Md and $file2Md"
